$wb = $excel.ActiveWorkbook

# --- Rename the five sheets to their new short names ---------------------
$wb.Worksheets.Item(1).Name = "NI03"
$wb.Worksheets.Item(2).Name = "NI26"
$wb.Worksheets.Item(3).Name = "NI70"
$wb.Worksheets.Item(4).Name = "NI82"
$wb.Worksheets.Item(5).Name = "NI93"

# --- Make the first sheet the active tab (was the fifth sheet) -----------
$wb.Worksheets.Item(1).Activate()

# --- Defined names: refresh every name so its RefersTo formula is
#     re-serialised against the new sheet names (drops the now-unneeded
#     quoting around names that no longer contain spaces) ------------------
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    $n.RefersTo = $n.RefersTo
}

# The built-in Print_Titles name isn't refreshed by a sheet rename, so it
# needs to be pointed at the new sheet name explicitly.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -like "*Print_Titles*") {
        $n.RefersTo = "=NI93!`$1:`$5"
    }
}

# --- Header/footer margins: nudge the stored inches value from the
#     points-derived figure (36.85 pt) to the exact 1.3 cm figure ---------
for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.PageSetup.HeaderMargin = 36.850393700787386
    $ws.PageSetup.FooterMargin = 36.850393700787386
}

# --- Drop the trailing space in the last sheet's right footer ------------
$wb.Worksheets.Item(5).PageSetup.RightFooter = "&8Page &P of &N"
